# -----------------------------------------------------------------------
# Adds 12 new experimental data rows (225-236) for the FeCoNiCr0.5Al(x)
# "SHS" (thermite) alloy series, sourced from 10.1007/s11837-019-03678-3,
# as described in the commit message:
#   "extracted data from `10.1007/s11837-019-03678-3`"
# Each of the 4 alloy compositions (Al0.6 / Al0.8 / Al1.0 / Al1.2) gets
# three rows: compressive yield stress, compressive fracture stress, and
# compressive ductility.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Excel constant used below (xlCenter), so J's number format/alignment
# matches the centered "#,##0" style used for these new rows.
$xlCenter = -4108

# Column B
$ws.Cells.Item(225, 2).Value = "FeCoNiCr0.5Al0.6"
$ws.Cells.Item(226, 2).Value = "FeCoNiCr0.5Al0.8"
$ws.Cells.Item(227, 2).Value = "FeCoNiCr0.5Al1.0"
$ws.Cells.Item(228, 2).Value = "FeCoNiCr0.5Al1.2"
$ws.Cells.Item(229, 2).Value = "FeCoNiCr0.5Al0.6"
$ws.Cells.Item(230, 2).Value = "FeCoNiCr0.5Al0.8"
$ws.Cells.Item(231, 2).Value = "FeCoNiCr0.5Al1.0"
$ws.Cells.Item(232, 2).Value = "FeCoNiCr0.5Al1.2"
$ws.Cells.Item(233, 2).Value = "FeCoNiCr0.5Al0.6"
$ws.Cells.Item(234, 2).Value = "FeCoNiCr0.5Al0.8"
$ws.Cells.Item(235, 2).Value = "FeCoNiCr0.5Al1.0"
$ws.Cells.Item(236, 2).Value = "FeCoNiCr0.5Al1.2"

# Column C
$ws.Cells.Item(225, 3).Value = "FCC+BCC"
$ws.Cells.Item(226, 3).Value = "FCC+BCC"
$ws.Cells.Item(227, 3).Value = "FCC+BCC+B2"
$ws.Cells.Item(228, 3).Value = "BCC+B2"
$ws.Cells.Item(229, 3).Value = "FCC+BCC"
$ws.Cells.Item(230, 3).Value = "FCC+BCC"
$ws.Cells.Item(231, 3).Value = "FCC+BCC+B2"
$ws.Cells.Item(232, 3).Value = "BCC+B2"
$ws.Cells.Item(233, 3).Value = "FCC+BCC"
$ws.Cells.Item(234, 3).Value = "FCC+BCC"
$ws.Cells.Item(235, 3).Value = "FCC+BCC+B2"
$ws.Cells.Item(236, 3).Value = "BCC+B2"

# Column D
$ws.Cells.Item(225, 4).Value = "SHS"
$ws.Cells.Item(226, 4).Value = "SHS"
$ws.Cells.Item(227, 4).Value = "SHS"
$ws.Cells.Item(228, 4).Value = "SHS"
$ws.Cells.Item(229, 4).Value = "SHS"
$ws.Cells.Item(230, 4).Value = "SHS"
$ws.Cells.Item(231, 4).Value = "SHS"
$ws.Cells.Item(232, 4).Value = "SHS"
$ws.Cells.Item(233, 4).Value = "SHS"
$ws.Cells.Item(234, 4).Value = "SHS"
$ws.Cells.Item(235, 4).Value = "SHS"
$ws.Cells.Item(236, 4).Value = "SHS"

# Column E
$ws.Cells.Item(225, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(226, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(227, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(228, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(229, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(230, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(231, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(232, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(233, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(234, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(235, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."
$ws.Cells.Item(236, 5).Value = "SHS thermite reaction of Co Cr Fe Ni oxides with Al."

# Column F
$ws.Cells.Item(225, 6).Value = "compressive yield stress"
$ws.Cells.Item(226, 6).Value = "compressive yield stress"
$ws.Cells.Item(227, 6).Value = "compressive yield stress"
$ws.Cells.Item(228, 6).Value = "compressive yield stress"
$ws.Cells.Item(229, 6).Value = "compressive fracture stress"
$ws.Cells.Item(230, 6).Value = "compressive fracture stress"
$ws.Cells.Item(231, 6).Value = "compressive fracture stress"
$ws.Cells.Item(232, 6).Value = "compressive fracture stress"
$ws.Cells.Item(233, 6).Value = "compressive ductility"
$ws.Cells.Item(234, 6).Value = "compressive ductility"
$ws.Cells.Item(235, 6).Value = "compressive ductility"
$ws.Cells.Item(236, 6).Value = "compressive ductility"

# Column G
$ws.Cells.Item(225, 7).Value = "EXP"
$ws.Cells.Item(226, 7).Value = "EXP"
$ws.Cells.Item(227, 7).Value = "EXP"
$ws.Cells.Item(228, 7).Value = "EXP"
$ws.Cells.Item(229, 7).Value = "EXP"
$ws.Cells.Item(230, 7).Value = "EXP"
$ws.Cells.Item(231, 7).Value = "EXP"
$ws.Cells.Item(232, 7).Value = "EXP"
$ws.Cells.Item(233, 7).Value = "EXP"
$ws.Cells.Item(234, 7).Value = "EXP"
$ws.Cells.Item(235, 7).Value = "EXP"
$ws.Cells.Item(236, 7).Value = "EXP"

# Column I
$ws.Cells.Item(225, 9).Value = 298
$ws.Cells.Item(226, 9).Value = 298
$ws.Cells.Item(227, 9).Value = 298
$ws.Cells.Item(228, 9).Value = 298
$ws.Cells.Item(229, 9).Value = 298
$ws.Cells.Item(230, 9).Value = 298
$ws.Cells.Item(231, 9).Value = 298
$ws.Cells.Item(232, 9).Value = 298
$ws.Cells.Item(233, 9).Value = 298
$ws.Cells.Item(234, 9).Value = 298
$ws.Cells.Item(235, 9).Value = 298
$ws.Cells.Item(236, 9).Value = 298

# Column J (value + number format/alignment to match style 42)
$ws.Cells.Item(225, 10).Value = 412000000
$ws.Cells.Item(225, 10).NumberFormat = "#,##0"
$ws.Cells.Item(225, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(226, 10).Value = 1214000000
$ws.Cells.Item(226, 10).NumberFormat = "#,##0"
$ws.Cells.Item(226, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(227, 10).Value = 1140000000
$ws.Cells.Item(227, 10).NumberFormat = "#,##0"
$ws.Cells.Item(227, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(228, 10).Value = 1000
$ws.Cells.Item(228, 10).NumberFormat = "#,##0"
$ws.Cells.Item(228, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(229, 10).Value = 2112000000
$ws.Cells.Item(229, 10).NumberFormat = "#,##0"
$ws.Cells.Item(229, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(230, 10).Value = 2923000000
$ws.Cells.Item(230, 10).NumberFormat = "#,##0"
$ws.Cells.Item(230, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(231, 10).Value = 2873000000
$ws.Cells.Item(231, 10).NumberFormat = "#,##0"
$ws.Cells.Item(231, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(232, 10).Value = 1832000000
$ws.Cells.Item(232, 10).NumberFormat = "#,##0"
$ws.Cells.Item(232, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(233, 10).Value = 42
$ws.Cells.Item(233, 10).NumberFormat = "#,##0"
$ws.Cells.Item(233, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(234, 10).Value = 37
$ws.Cells.Item(234, 10).NumberFormat = "#,##0"
$ws.Cells.Item(234, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(235, 10).Value = 31
$ws.Cells.Item(235, 10).NumberFormat = "#,##0"
$ws.Cells.Item(235, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(236, 10).Value = 16
$ws.Cells.Item(236, 10).NumberFormat = "#,##0"
$ws.Cells.Item(236, 10).HorizontalAlignment = $xlCenter

# Column K (only row 228 has a value)
$ws.Cells.Item(228, 11).Value = 100

# Column L
$ws.Cells.Item(225, 12).Value = "Pa"
$ws.Cells.Item(226, 12).Value = "Pa"
$ws.Cells.Item(227, 12).Value = "Pa"
$ws.Cells.Item(228, 12).Value = "Pa"
$ws.Cells.Item(229, 12).Value = "Pa"
$ws.Cells.Item(230, 12).Value = "Pa"
$ws.Cells.Item(231, 12).Value = "Pa"
$ws.Cells.Item(232, 12).Value = "Pa"
$ws.Cells.Item(233, 12).Value = "%"
$ws.Cells.Item(234, 12).Value = "%"
$ws.Cells.Item(235, 12).Value = "%"
$ws.Cells.Item(236, 12).Value = "%"

# Column M
$ws.Cells.Item(225, 13).Value = "P3461"
$ws.Cells.Item(226, 13).Value = "P3461"
$ws.Cells.Item(227, 13).Value = "P3461"
$ws.Cells.Item(228, 13).Value = "F2"
$ws.Cells.Item(229, 13).Value = "P3461"
$ws.Cells.Item(230, 13).Value = "P3461"
$ws.Cells.Item(231, 13).Value = "P3461"
$ws.Cells.Item(232, 13).Value = "P3461"
$ws.Cells.Item(233, 13).Value = "P3461"
$ws.Cells.Item(234, 13).Value = "P3461"
$ws.Cells.Item(235, 13).Value = "P3461"
$ws.Cells.Item(236, 13).Value = "P3461"

# Column N
$ws.Cells.Item(225, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(226, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(227, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(228, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(229, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(230, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(231, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(232, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(233, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(234, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(235, 14).Value = "10.1007/s11837-019-03678-3"
$ws.Cells.Item(236, 14).Value = "10.1007/s11837-019-03678-3"

# Update the sheet view to reflect where the user was working (matches the
# commit's sheetView/selection change): scroll so column G / row 212 is the
# top-left visible cell, with N240 as the active selection.
$excel.ActiveWindow.ScrollRow = 212
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("N240").Select()

# Best-effort: also nudge the workbook window position to match the
# commit's bookViews/workbookView xWindow/yWindow (cosmetic only).
$excel.ActiveWindow.Left = 9800
$excel.ActiveWindow.Top = 2220
